$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (CBO decoherence envelope) - 9d column (D): 28.7 -> 9.3
$ws.Range("D12").Value = 9.3000000000000007

# Row 14 (Fixed k loss) - HighKick column (C): 39.1 -> 4.9
$ws.Range("C14").Value = 4.9000000000000004

# Row 14 (Fixed k loss) - Endgame column (E): 28.9 -> 0.1
$ws.Range("E14").Value = 0.1

# Update the active selection on the sheet view to F18
$ws.Range("F18").Select()
